# Add calc U I angles
# Insert a new column before column C ("Ua"), shifting all existing
# per-point measurement columns (Ua..MTE γ изм) one column to the right,
# and use the freed-up column C for a new "freq" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting at column C pushes the existing C:AI block to D:AJ, carrying
# over cell formatting (matches the widened dimension/cols seen in the
# target workbook).
$ws.Columns("C").Insert() | Out-Null

# Populate the header for the newly freed column C.
$ws.Range("C1").Value = "freq"

# Leave the new header cell selected, as in the saved workbook.
$ws.Range("C1").Select() | Out-Null
